$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-All "INTERNISTA" "CHIRURG"
Replace-All "lightskyblue" "yellow"
Replace-All "TEST GAWEŁ" "TEST TEST"
Replace-All "08.07.2023 13:07" "16.07.2023 17:41"
Replace-All "Wituś 1" "1"
Replace-All "IDOCUMENTSERVICE: DOCUMENTSERVICE" "TEST"
Replace-All "13:07" "17:41"
Replace-All "Badanie eratury ciała" "Badanie 36eratury ciała"

# Fill in the blank vital-sign value cells of the "KARTA OBSERWACJI" table
# (table #4): the value cell is one column to the right of the row label.
$vitals = $d.Tables.Item(4)
$valueByRow = @{
    3  = "120"  # Cisnienie tetnicze skurczowe
    4  = "80"   # Cisnienie tetnicze rozkurczowe
    5  = "80"   # Tetno
    17 = "99"   # Saturacja O2
    21 = "15"   # Wedlug GCS
    35 = "36"   # Temperatura
}
foreach ($row in $valueByRow.Keys) {
    $cell = $vitals.Cell($row, 5)
    $cell.Range.Text = $valueByRow[$row]
}
